$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "ZUpna1219AS1"
$ws.Range("C4").Value = "PNABU-L3-ZU-022"
$ws.Range("D4").Value = "ZUpna-1219AS-1"
$ws.Range("H4").Value = "25 Nov 2023"
